$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.700.30'
$ws.Range("E2").Value = '  +0.79%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.657.20'
$ws.Range("E3").Value = '  +1.06%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.000'
$ws.Range("E5").Value = '  +0.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.46'
$ws.Range("E6").Value = '  -0.27%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3815'
$ws.Range("E7").Value = '  +0.49%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3614'
$ws.Range("E8").Value = '  -0.48%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '51.24'
$ws.Range("E9").Value = '  -0.85%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08204'
$ws.Range("E10").Value = '  +0.12%  '

$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.239'
$ws.Range("E11").Value = '  +0.25%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  +0.11%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.66'
$ws.Range("E13").Value = '  +0.69%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.517'
$ws.Range("E14").Value = '  +0.68%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.428'
$ws.Range("E15").Value = '  +0.32%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001235'
$ws.Range("E16").Value = '  -0.56%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.652.88'
$ws.Range("E17").Value = '  +1.26%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '97.53'
$ws.Range("E18").Value = '  +2.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06998'
$ws.Range("E19").Value = '  +0.88%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.832'
$ws.Range("E20").Value = '  +3.54%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.74'
$ws.Range("E21").Value = '  +1.28%  '

$ws.Range("E22").Value = '  +0.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.83'
$ws.Range("E23").Value = '  +2.30%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.709.47'
$ws.Range("E24").Value = '  +0.80%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.514'
$ws.Range("E25").Value = '  +0.72%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.056'
$ws.Range("E26").Value = '  -0.38%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.30'
$ws.Range("E27").Value = '  +0.73%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '153.20'
$ws.Range("E28").Value = '  +1.09%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.213'
$ws.Range("E29").Value = '  -1.05%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.78'
$ws.Range("E30").Value = '  +1.08%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.837.09'
$ws.Range("E31").Value = '  +1.18%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.986'
$ws.Range("E32").Value = '  +4.76%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.199'
$ws.Range("E33").Value = '  +0.45%  '

$ws.Range("E34").Value = '  +5.81%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.062'
$ws.Range("E35").Value = '  +0.10%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02818'
$ws.Range("E36").Value = '  +1.80%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2526'
$ws.Range("E37").Value = '  +0.95%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.130'
$ws.Range("E38").Value = '  +1.61%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.08795'
$ws.Range("E39").Value = '  +0.13%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.07091'
$ws.Range("E40").Value = '  -0.43%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.11'
$ws.Range("E41").Value = '  +7.76%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7041'
$ws.Range("E42").Value = '  -0.18%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.337'
$ws.Range("E43").Value = '  -0.10%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.10'
$ws.Range("E44").Value = '  +1.73%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6545'
$ws.Range("E45").Value = '  +0.00%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.324'
$ws.Range("E46").Value = '  +1.54%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9994'
$ws.Range("E47").Value = '  +0.03%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.973'
$ws.Range("E48").Value = '  +0.15%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07953'
$ws.Range("E49").Value = '  -0.33%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '128.13'
$ws.Range("E50").Value = '  -0.42%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.186'
$ws.Range("E51").Value = '  -0.34%  '
